# Fixed test template version: remove the orange "M2Doc version mismatch"
# warning runs (and surrounding padding spaces) from the first paragraph,
# leaving only the introductory sentence.
$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
